$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto pricing/volume figures (and the BKEXToken/KickToken row swap)
# exactly as scraped by the symbol-list GitHub Action. Cells hold text-formatted
# numbers/percentages, so write with a leading apostrophe to keep them as text
# (matching the original inline-string cell type) and reset Style so no stray
# number-format is left behind.
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "256.39"
Set-TextCell "E2" "0.43%"
Set-TextCell "D3" "26.97"
Set-TextCell "E3" "-4.12%"
Set-TextCell "D4" "4.723"
Set-TextCell "E4" "-10.02%"
Set-TextCell "D5" "0.05937"
Set-TextCell "E5" "1.49%"
Set-TextCell "D6" "6.660"
Set-TextCell "E6" "-0.68%"
Set-TextCell "D7" "0.8680"
Set-TextCell "E7" "-0.02%"
Set-TextCell "D8" "0.9537"
Set-TextCell "E8" "-7.82%"
Set-TextCell "D9" "0.1404"
Set-TextCell "E9" "-0.48%"
Set-TextCell "D10" "0.03959"
Set-TextCell "E10" "14.11%"
Set-TextCell "D11" "0.07167"
Set-TextCell "E11" "0.51%"
Set-TextCell "D12" "0.03186"
Set-TextCell "E12" "-0.07%"
Set-TextCell "D13" "0.09254"
Set-TextCell "E13" "0.30%"
Set-TextCell "D14" "0.001544"
Set-TextCell "E14" "0.19%"
Set-TextCell "D15" "0.0006071"
Set-TextCell "E15" "-94.28%"
Set-TextCell "D16" "0.006076"
Set-TextCell "E16" "4.26%"
Set-TextCell "D17" "3.483"
Set-TextCell "E17" "-0.42%"
Set-TextCell "D18" "3.200"
Set-TextCell "E18" "-1.03%"
Set-TextCell "E19" "-0.19%"
Set-TextCell "E20" "-1.51%"
Set-TextCell "E21" "-1.39%"
Set-TextCell "D22" "3.804"
Set-TextCell "E22" "7.17%"
Set-TextCell "D23" "0.04221"
Set-TextCell "E23" "1.29%"
Set-TextCell "E24" "2.42%"
Set-TextCell "D25" "0.001220"
Set-TextCell "E25" "-0.28%"
Set-TextCell "D26" "0.004496"
Set-TextCell "E26" "-7.87%"
Set-TextCell "E27" "-0.04%"
Set-TextCell "D28" "0.0001937"
Set-TextCell "E28" "-0.07%"
Set-TextCell "D40" "0.03823"
Set-TextCell "E40" "0.13%"
Set-TextCell "B41" "KickToken"
Set-TextCell "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D41" "0.005809"
Set-TextCell "E41" "1.35%"
Set-TextCell "B42" "BKEXToken"
Set-TextCell "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D42" "0.1100"
Set-TextCell "E42" "-0.16%"
Set-TextCell "D43" "0.002105"
Set-TextCell "E43" "-9.99%"
Set-TextCell "D44" "0.01056"
Set-TextCell "E44" "9.12%"
Set-TextCell "D45" "0.00005501"
Set-TextCell "E45" "5.16%"
Set-TextCell "E46" "-0.03%"
Set-TextCell "D47" "0.08853"
Set-TextCell "E47" "-4.83%"
Set-TextCell "D48" "0.002391"
Set-TextCell "E48" "11.11%"
Set-TextCell "D49" "0.00002100"
Set-TextCell "E49" "-0.03%"
Set-TextCell "D50" "0.0002000"
Set-TextCell "E50" "-0.03%"
